$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------

# "N°1" -> "Tout" : this contribution line now applies to the whole feature
$ws.Range("E5").Value = "Tout"

# New contribution line (row 11): Anthony, 08/01/2021, feature #7, "Tout"
$ws.Range("B11").Value = 44204
$ws.Range("C11").Value = "Anthony"
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = "Tout"

# Column B carries a date number format on every other row; copy it down
# onto the new row instead of leaving the plain column default.
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Conditional formatting -------------------------------------------------
# The "Nom" column is colour-coded per author. Re-apply the colour rules so
# the highlighting also covers the newly inserted row (public list display).

$colLucas    = 49407      # RGB(255,192,0)   -> FFC000
$colThomas   = 16711935   # RGB(255,0,255)   -> FF00FF
$colAnthony  = 10498160   # RGB(112,48,160)  -> 7030A0
$colFrancois = 12611584   # RGB(0,112,192)   -> 0070C0

$bigRange = $ws.Range("C1:C10,C17:C1048576")
$fc = $bigRange.FormatConditions.Add(1, 3, '"Lucas"')
$fc.Interior.Color = $colLucas
$fc = $bigRange.FormatConditions.Add(1, 3, '"Thomas"')
$fc.Interior.Color = $colThomas
$fc = $bigRange.FormatConditions.Add(1, 3, '"Anthony"')
$fc.Interior.Color = $colAnthony
$fc = $bigRange.FormatConditions.Add(1, 3, '"François"')
$fc.Interior.Color = $colFrancois

$newRow = $ws.Range("C11")
$fc = $newRow.FormatConditions.Add(1, 3, '"Lucas"')
$fc.Interior.Color = $colLucas
$fc = $newRow.FormatConditions.Add(1, 3, '"Thomas"')
$fc.Interior.Color = $colThomas
$fc = $newRow.FormatConditions.Add(1, 3, '"Anthony"')
$fc.Interior.Color = $colAnthony
$fc = $newRow.FormatConditions.Add(1, 3, '"François"')
$fc.Interior.Color = $colFrancois

# --- View bookkeeping -------------------------------------------------------
$ws.Range("B14").Select() | Out-Null
